$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3190
$ws.Range("C2").Value = 75.43000000000001

$ws.Range("B3").Value = 623
$ws.Range("C3").Value = 14.73

$ws.Range("B4").Value = 260
$ws.Range("C4").Value = 6.15

$ws.Range("B5").Value = 107
$ws.Range("C5").Value = 2.53

$ws.Range("B6").Value = 49
$ws.Range("C6").Value = 1.16
